$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 158.3
$ws.Range("J33").Value = 117.5
$ws.Range("L33").Value = 117.5
$ws.Range("N33").Value = -575.5
$ws.Range("H70").Value = 4981.757
$ws.Range("I70").Value = 2221.9697
$ws.Range("J70").Value = 27750
$ws.Range("K70").Value = 6665.909100000001
$ws.Range("L70").Value = 83250
$ws.Range("M70").Value = -6395.909100000001
$ws.Range("N70").Value = -83790
$ws.Range("H73").Value = 4981.757
$ws.Range("I73").Value = 2221.9697
$ws.Range("J73").Value = 27750
$ws.Range("K73").Value = 6665.909100000001
$ws.Range("L73").Value = 83250
$ws.Range("M73").Value = -5729.909100000001
$ws.Range("N73").Value = -85122
$ws.Range("H86").Value = 2798.8333
$ws.Range("I86").Value = 2958.6
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2958.6
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -1835.6
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2798.8333
$ws.Range("I89").Value = 2958.6
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 14793
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -9177
$ws.Range("N89").Value = -21232
$ws.Range("H112").Value = 2026
$ws.Range("J112").Value = 2027.1818
$ws.Range("L112").Value = 6081.5454
$ws.Range("N112").Value = -8297.545399999999
$ws.Range("H121").Value = 796
$ws.Range("J121").Value = 945
$ws.Range("L121").Value = 2835
$ws.Range("N121").Value = -6329
$ws.Range("H132").Value = 1344.6586
$ws.Range("I132").Value = 1214.1578
$ws.Range("J132").Value = 2997.6667
$ws.Range("K132").Value = 3642.4734
$ws.Range("L132").Value = 8993.000100000001
$ws.Range("M132").Value = -1112.4734
$ws.Range("N132").Value = -14053.0001
$ws.Range("H138").Value = 3518.4856
$ws.Range("J138").Value = 3468.3157
$ws.Range("L138").Value = 10404.9471
$ws.Range("N138").Value = -20684.9471
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1112696.2
$ws.Range("I2").Value = 2778777.5
$ws.Range("J2").Value = 1975.3334
$ws.Range("K2").Value = 2778777.5
$ws.Range("L2").Value = 1975.3334
$ws.Range("M2").Value = -2778664.5
$ws.Range("N2").Value = -2201.3334
$ws.Range("H32").Value = 3900.7834
$ws.Range("I32").Value = 3100.963
$ws.Range("K32").Value = 3100.963
$ws.Range("M32").Value = -2813.963
$ws.Range("H45").Value = 1815.5385
$ws.Range("I45").Value = 1355.5
$ws.Range("K45").Value = 1355.5
$ws.Range("M45").Value = -978.5
$ws.Range("H61").Value = 2627.2666
$ws.Range("I61").Value = 1425.8334
$ws.Range("J61").Value = 7433
$ws.Range("K61").Value = 1425.8334
$ws.Range("L61").Value = 7433
$ws.Range("M61").Value = -1213.8334
$ws.Range("N61").Value = -7857
$ws.Range("H74").Value = 1859.4839
$ws.Range("I74").Value = 1502.7391
$ws.Range("K74").Value = 1502.7391
$ws.Range("M74").Value = -628.7391
$ws.Range("H77").Value = 1859.4839
$ws.Range("I77").Value = 1502.7391
$ws.Range("K77").Value = 7513.6955
$ws.Range("M77").Value = -3145.6955
$ws.Range("H116").Value = 1112696.2
$ws.Range("I116").Value = 2778777.5
$ws.Range("J116").Value = 1975.3334
$ws.Range("K116").Value = 2778777.5
$ws.Range("L116").Value = 1975.3334
$ws.Range("M116").Value = -2776483.5
$ws.Range("N116").Value = -6563.3334
$ws.Range("H132").Value = 1505.3214
$ws.Range("I132").Value = 1216.5652
$ws.Range("J132").Value = 2833.6
$ws.Range("K132").Value = 3649.6956
$ws.Range("L132").Value = 8500.799999999999
$ws.Range("M132").Value = -1119.6956
$ws.Range("N132").Value = -13560.8
$ws.Range("H136").Value = 2627.2666
$ws.Range("I136").Value = 1425.8334
$ws.Range("J136").Value = 7433
$ws.Range("K136").Value = 4277.5002
$ws.Range("L136").Value = 22299
$ws.Range("M136").Value = -1727.5002
$ws.Range("N136").Value = -27399
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1112696.2
$ws.Range("I3").Value = 2778777.5
$ws.Range("J3").Value = 1975.3334
$ws.Range("K3").Value = 2778777.5
$ws.Range("L3").Value = 1975.3334
$ws.Range("M3").Value = -2778663.5
$ws.Range("N3").Value = -2203.3334
$ws.Range("H105").Value = 2233.8
$ws.Range("I105").Value = 2229.0715
$ws.Range("K105").Value = 2229.0715
$ws.Range("M105").Value = -482.0715
$ws.Range("H134").Value = 7831.222
$ws.Range("I134").Value = 7831.222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 23493.666
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -20958.666
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H31").Value = 2215.6956
$ws.Range("I31").Value = 1955.9231
$ws.Range("J31").Value = 2553.4
$ws.Range("K31").Value = 1955.9231
$ws.Range("L31").Value = 2553.4
$ws.Range("M31").Value = -1660.9231
$ws.Range("N31").Value = -3143.4
$ws.Range("H34").Value = 2215.6956
$ws.Range("I34").Value = 1955.9231
$ws.Range("J34").Value = 2553.4
$ws.Range("K34").Value = 1955.9231
$ws.Range("L34").Value = 2553.4
$ws.Range("M34").Value = -1753.9231
$ws.Range("N34").Value = -2957.4
$ws.Range("H99").Value = 2959.6667
$ws.Range("I99").Value = 1827.4
$ws.Range("K99").Value = 1827.4
$ws.Range("M99").Value = -329.4000000000001
$ws.Range("H105").Value = 1124.5
$ws.Range("I105").Value = 1249.75
$ws.Range("J105").Value = 874
$ws.Range("K105").Value = 1249.75
$ws.Range("L105").Value = 874
$ws.Range("M105").Value = 497.25
$ws.Range("N105").Value = -4368
$ws.Range("H126").Value = 2959.6667
$ws.Range("I126").Value = 1827.4
$ws.Range("K126").Value = 5482.200000000001
$ws.Range("M126").Value = -3012.200000000001
$ws.Range("H132").Value = 2547.5483
$ws.Range("I132").Value = 1843.8096
$ws.Range("J132").Value = 4025.4
$ws.Range("K132").Value = 5531.4288
$ws.Range("L132").Value = 12076.2
$ws.Range("M132").Value = -3001.4288
$ws.Range("N132").Value = -17136.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1630.1111
$ws.Range("J118").Value = 1973.6666
$ws.Range("L118").Value = 5920.9998
$ws.Range("N118").Value = -8406.9998
$ws.Range("H131").Value = 20118.973
$ws.Range("J131").Value = 22544.469
$ws.Range("L131").Value = 67633.40700000001
$ws.Range("N131").Value = -77713.40700000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1905.5333
$ws.Range("J97").Value = 1924.25
$ws.Range("L97").Value = 1924.25
$ws.Range("N97").Value = -2916.25
$ws.Range("H132").Value = 1674839.4
$ws.Range("I132").Value = 2749252
$ws.Range("K132").Value = 8247756
$ws.Range("M132").Value = -8245226
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1810.7142
$ws.Range("I46").Value = 1096.125
$ws.Range("J46").Value = 2763.5
$ws.Range("K46").Value = 1096.125
$ws.Range("L46").Value = 2763.5
$ws.Range("M46").Value = -908.125
$ws.Range("N46").Value = -3139.5
$ws.Range("H55").Value = 776.2
$ws.Range("I55").Value = 800.5
$ws.Range("J55").Value = 760
$ws.Range("K55").Value = 800.5
$ws.Range("L55").Value = 760
$ws.Range("M55").Value = -627.5
$ws.Range("N55").Value = -1106
$ws.Range("H122").Value = 4824.643
$ws.Range("I122").Value = 4516.5557
$ws.Range("K122").Value = 13549.6671
$ws.Range("M122").Value = -11099.6671
$ws.Range("H132").Value = 4162
$ws.Range("I132").Value = 4016.4443
$ws.Range("K132").Value = 12049.3329
$ws.Range("M132").Value = -9519.332900000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14330
$ws.Range("J41").Value = 14330
$ws.Range("L41").Value = 14330
$ws.Range("N41").Value = -15110
$ws.Range("H108").Value = 24000
$ws.Range("J108").Value = 24000
$ws.Range("L108").Value = 24000
$ws.Range("N108").Value = -31680
$ws.Range("H109").Value = 78318.75
$ws.Range("J109").Value = 78318.75
$ws.Range("L109").Value = 78318.75
$ws.Range("N109").Value = -81092.75
$ws.Range("H132").Value = 1922.7
$ws.Range("I132").Value = 1183.9
$ws.Range("J132").Value = 3400.3
$ws.Range("K132").Value = 3551.7
$ws.Range("L132").Value = 10200.9
$ws.Range("M132").Value = -1021.7
$ws.Range("N132").Value = -15260.9
$ws.Range("H136").Value = 10289290
$ws.Range("I136").Value = 16341087
$ws.Range("J136").Value = 1235.4
$ws.Range("K136").Value = 49023261
$ws.Range("L136").Value = 3706.2
$ws.Range("M136").Value = -49020711
$ws.Range("N136").Value = -8806.200000000001
